$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Grants): "None Found" -> specific grant numbers (kept as text,
# not auto-converted to a number). We temporarily format the cell as Text so
# the numeric-looking string "1419282" is stored as a shared string, then
# reset the style back to Normal so no style id is left on the cell.
$grantsRows24 = @(2,3,4)
foreach ($r in $grantsRows24) {
    $c = $ws.Range("B$r")
    $c.NumberFormat = "@"
    $c.Value = "1419282"
    $c.Style = "Normal"
}

$grantsRows56 = @(5,6)
foreach ($r in $grantsRows56) {
    $c = $ws.Range("B$r")
    $c.NumberFormat = "@"
    $c.Value = "P42ES007380, R03OD030603, 1419282, 2020026"
    $c.Style = "Normal"
}

# Column N (Publication Year): was stored as text "2020"/"2021", now a real
# number.
$ws.Range("N2").Value = 2020
$ws.Range("N3").Value = 2020
$ws.Range("N4").Value = 2020
$ws.Range("N5").Value = 2021
$ws.Range("N6").Value = 2021

# Column V (Comparison): data removed entirely for all data rows.
$ws.Range("V2").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("V6").ClearContents()
